$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawl was re-run / re-saved later the same day: every row's
# "timestamp" column (O2:O530) moves from 06:52:16 to 12:59:19.
$ws.Range("O2:O530").Value = "2023-01-03 12:59:19"

# A handful of products came back in stock between the two crawls, so
# their aria-label (column M) no longer carries the
# "- Online kein Bestand" (out of stock online) suffix.
$ws.Range("M324").Value = "Fairtrade Papaya 1 Stück 2.85 Schweizer Franken"
$ws.Range("M360").Value = "Naturaplan Bio Shiitake-Pilze ca. 100g 3.20 Schweizer Franken"
$ws.Range("M430").Value = "Birnen Rocha IP-Suisse ca. 1kg 4.50 Schweizer Franken"
